# Update the dSF (column F) values on the active sheet to reflect the
# repulled data / recalculated mean values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = 8
    5  = -2
    6  = 1
    7  = -3
    8  = 1
    9  = -1
    10 = 1
    14 = -1
    15 = 7
    17 = 2
    18 = 3
    19 = -2
    20 = -5
    21 = 3
    23 = 5
    24 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
